$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (FAPs / Tac1 / Tacr1 / ECs) values
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.213143333333333
$ws.Range("H2").Value = 3.63943
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 0.869364
$ws.Range("N2").Value = 1.738728
$ws.Range("O2").Value = 0.8481925485270655
$ws.Range("P2").Value = 0.788353574987826
$ws.Range("Q2").Value = 1.05466314084
$ws.Range("R2").Value = 6.327978845040001
$ws.Range("S2").Value = 0.8481925485270655
$ws.Range("T2").Value = 0.788353574987826

# Add new row 3 (FAPs / Tac1 / Tacr1 / Neutro)
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Tac1"
$ws.Range("C3").Value = "Tacr1"
$ws.Range("D3").Value = "Neutro"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.213143333333333
$ws.Range("H3").Value = 3.63943
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1555966666666667
$ws.Range("N3").Value = 0.46679
$ws.Range("O3").Value = 0.1518074514729346
$ws.Range("P3").Value = 0.211646425012174
$ws.Range("Q3").Value = 0.1887610588555556
$ws.Range("R3").Value = 1.6988495297
$ws.Range("S3").Value = 0.1518074514729346
$ws.Range("T3").Value = 0.211646425012174
